$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Helper: write literal text into a cell without letting the engine coerce
# numeric-looking strings (like "1.") into numbers. We build the text via a
# formula in a scratch cell, copy it, and paste-special as values into the
# destination -- this yields a plain shared-string cell with no stray style.
function Set-LiteralText($range, [string]$text) {
    $scratch = $ws.Range("Z100")
    $scratch.Formula = "=""" + $text + """"
    $scratch.Copy()
    $range.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
    $scratch.ClearContents()
}

# ---- Table 1 (rows 2-5) ----
$ws.Range("B2").Value = "I"
$ws.Range("C2").Value = "II"
$ws.Range("D2").Value = "III"
$ws.Range("E2").Value = "IV"
$ws.Range("A3").Value = "a"
$ws.Range("C3").Value = 22
$ws.Range("A4").Value = "b"
$ws.Range("C4").Value = 66
$ws.Range("A5").Value = "c"
$ws.Range("C5").Value = 4

# ---- Table 2 (rows 10-13) ----
$ws.Range("B10").Value = "I"
$ws.Range("C10").Value = "II"
$ws.Range("D10").Value = "III"
$ws.Range("E10").Value = "IV"
$ws.Range("A11").Value = "a"
$ws.Range("C11").Value = 51
$ws.Range("A12").Value = "b"
$ws.Range("C12").Value = 443
$ws.Range("A13").Value = "c"
$ws.Range("C13").Value = 5

# ---- Table 3 (rows 18-21) ----
$ws.Range("B18").Value = "I"
$ws.Range("C18").Value = "II"
$ws.Range("D18").Value = "III"
$ws.Range("E18").Value = "IV"
$ws.Range("A19").Value = "a"
$ws.Range("C19").Value = 37
$ws.Range("A20").Value = "b"
$ws.Range("C20").Value = 3462
$ws.Range("A21").Value = "c"
$ws.Range("C21").Value = 51

# ---- Section captions (added last so they land at the tail of the shared
# string table, matching the source order) ----
Set-LiteralText $ws.Range("A1") "1."
Set-LiteralText $ws.Range("A9") "2."
Set-LiteralText $ws.Range("A17") "3."

# ---- Turn the three ranges into real Excel Tables ----
$tbl1 = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $ws.Range("B2:E5"), $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl2 = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $ws.Range("B10:E13"), $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl3 = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $ws.Range("B18:E21"), $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)

# Rename starting from the last table so none of the still-default-named
# tables collide with a name we are about to assign.
$tbl3.Name = "Table245"
$tbl2.Name = "Table24"
$tbl1.Name = "Table2"

# ---- Column widths A:F ----
$ws.Range("A1:F1").EntireColumn.ColumnWidth = 10.2

# ---- Final selection ----
$ws.Range("C22").Select()
